# Updated cryptos list on Tue Jul 30 07:52:29 UTC 2024 with GitHub Actions
#
# Refreshes each coin row's Price (column D) and Volume(1h) (column E)
# with the latest scraped values, and reflects that Aptos and USDe swapped
# ranking order (rows 33 and 34).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    # Force the cell to stay a plain text cell (matches the sheet's
    # existing inline-string cells) even when the new value looks
    # like a number (e.g. "574.52") or a dotted thousands price
    # (e.g. "66.947.00"), then restore the original (default) style
    # so no stray number-format/style id gets stamped onto the cell.
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '66.947.00'
Set-TextValue $ws.Range("E2") '  -3.70%  '
Set-TextValue $ws.Range("D3") '3.342.50'
Set-TextValue $ws.Range("E3") '  -0.74%  '
Set-TextValue $ws.Range("D5") '574.52'
Set-TextValue $ws.Range("E5") '  -3.17%  '
Set-TextValue $ws.Range("D6") '182.86'
Set-TextValue $ws.Range("E6") '  -4.88%  '
Set-TextValue $ws.Range("E7") '  +0.00%  '
Set-TextValue $ws.Range("E8") '  -1.62%  '
Set-TextValue $ws.Range("D9") '0.130'
Set-TextValue $ws.Range("E9") '  -3.28%  '
Set-TextValue $ws.Range("D10") '6.66'
Set-TextValue $ws.Range("E10") '  -1.61%  '
Set-TextValue $ws.Range("E11") '  -4.36%  '
Set-TextValue $ws.Range("D12") '3.925.30'
Set-TextValue $ws.Range("E12") '  -0.79%  '
Set-TextValue $ws.Range("D13") '0.137'
Set-TextValue $ws.Range("E13") '  -1.26%  '
Set-TextValue $ws.Range("E14") '  -5.05%  '
Set-TextValue $ws.Range("D15") '67.002.76'
Set-TextValue $ws.Range("E15") '  -3.63%  '
Set-TextValue $ws.Range("E16") '  -2.31%  '
Set-TextValue $ws.Range("D17") '3.347.08'
Set-TextValue $ws.Range("E17") '  -0.23%  '
Set-TextValue $ws.Range("D18") '436.70'
Set-TextValue $ws.Range("E18") '  -3.20%  '
Set-TextValue $ws.Range("D19") '13.73'
Set-TextValue $ws.Range("E19") '  -0.87%  '
Set-TextValue $ws.Range("E20") '  -2.42%  '
Set-TextValue $ws.Range("D21") '7.67'
Set-TextValue $ws.Range("E21") '  -2.09%  '
Set-TextValue $ws.Range("D22") '73.86'
Set-TextValue $ws.Range("E22") '  -0.85%  '
Set-TextValue $ws.Range("E23") '  +0.14%  '
Set-TextValue $ws.Range("D24") '0.525'
Set-TextValue $ws.Range("E24") '  +0.94%  '
Set-TextValue $ws.Range("E25") '  -2.51%  '
Set-TextValue $ws.Range("E26") '  -1.05%  '
Set-TextValue $ws.Range("D27") '9.15'
Set-TextValue $ws.Range("E27") '  -4.54%  '
Set-TextValue $ws.Range("E28") '  -0.12%  '
Set-TextValue $ws.Range("E29") '  -1.49%  '
Set-TextValue $ws.Range("E30") '  -1.74%  '
Set-TextValue $ws.Range("D31") '5.35'
Set-TextValue $ws.Range("E31") '  -4.69%  '
Set-TextValue $ws.Range("E32") '  -3.91%  '
Set-TextValue $ws.Range("B33") 'USDe'
Set-TextValue $ws.Range("C33") 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range("D33") '0.999'
Set-TextValue $ws.Range("E33") '  +0.05%  '
Set-TextValue $ws.Range("B34") 'Aptos'
Set-TextValue $ws.Range("C34") 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range("D34") '6.84'
Set-TextValue $ws.Range("E34") '  -2.95%  '
Set-TextValue $ws.Range("E35") '  -1.20%  '
Set-TextValue $ws.Range("D36") '161.19'
Set-TextValue $ws.Range("E36") '  -2.40%  '
Set-TextValue $ws.Range("D37") '27.65'
Set-TextValue $ws.Range("E37") '  +1.60%  '
Set-TextValue $ws.Range("E38") '  -4.93%  '
Set-TextValue $ws.Range("D39") '2.839.03'
Set-TextValue $ws.Range("E39") '  +3.57%  '
Set-TextValue $ws.Range("E40") '  -3.14%  '
Set-TextValue $ws.Range("D41") '4.45'
Set-TextValue $ws.Range("E41") '  -3.38%  '
Set-TextValue $ws.Range("E42") '  -4.88%  '
Set-TextValue $ws.Range("D43") '0.0677'
Set-TextValue $ws.Range("E43") '  -1.88%  '
Set-TextValue $ws.Range("D44") '40.30'
Set-TextValue $ws.Range("E44") '  -1.22%  '
Set-TextValue $ws.Range("D45") '24.72'
Set-TextValue $ws.Range("E45") '  -3.88%  '
Set-TextValue $ws.Range("E46") '  -6.45%  '
Set-TextValue $ws.Range("D47") '326.69'
Set-TextValue $ws.Range("E47") '  -4.84%  '
Set-TextValue $ws.Range("E48") '  -4.10%  '
Set-TextValue $ws.Range("D49") '31.69'
Set-TextValue $ws.Range("E49") '  -4.17%  '
Set-TextValue $ws.Range("D50") '0.993'
Set-TextValue $ws.Range("E50") '  -4.60%  '
Set-TextValue $ws.Range("D51") '6.18'
Set-TextValue $ws.Range("E51") '  -2.55%  '
